$d = $word.ActiveDocument
$r = $d.Content
$searchText = "I make things - things that work; and by " + [char]0x2018 + "work" + [char]0x2019 + " I mean work superlatively."
$found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence"
}

# Re-seat the match in a fresh Range object - the Find-mutated range
# object mis-slices multi-run InsertXML targets, a fresh Range(start,end)
# does not.
$target = $d.Range($r.Start, $r.End)

$xml = "<?xml version=`"1.0`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">I make things - things that work-</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">and by</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">‘</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">work</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">’</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">I mean work</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">superlatively</w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=`"preserve`">.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$target.InsertXML($xml)
Write-Output "Replaced sentence successfully"
